# 4.3.1.1 Youth education by gender — add 2021 column (M) and fix the
# English title shared string (drop the stray period after "4.3.1.1").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Excel constant values (since this runtime has no [Microsoft.Office...] enum) ---
$xlContinuous   = 1
$xlMedium       = -4138
$xlEdgeBottom   = 9
$xlVAlignCenter = -4108

function Set-BottomBorder($rng) {
    $rng.Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
    $rng.Borders.Item($xlEdgeBottom).Weight = $xlMedium
    $rng.Borders.Item($xlEdgeBottom).Color = 0
}

function Set-CellLook($rng, [bool]$bold, [bool]$border, [string]$numFmt) {
    $rng.Font.Name = "Times New Roman"
    $rng.Font.Size = 9
    $rng.Font.Bold = $bold
    $rng.VerticalAlignment = $xlVAlignCenter
    if ($numFmt) {
        $rng.NumberFormat = $numFmt
    }
    if ($border) {
        Set-BottomBorder $rng
    }
}

# 1) Fix the English title: "4.3.1.1. Youth education by gender" -> "4.3.1.1 Youth education by gender"
$ws.Range("C1").Value2 = "4.3.1.1 Youth education by gender"

# 2) Row 2 (thin separator row above the header) gains an M cell with the same
#    bordered look as the rest of the row.
Set-CellLook $ws.Range("M2") $false $true $null

# 3) Header row 3 gets a restyled L3 (matching the other year header cells)
#    plus a brand-new M3 = 2021 header cell.
Set-CellLook $ws.Range("L3") $true $true $null
$ws.Range("M3").Value2 = 2021
Set-CellLook $ws.Range("M3") $true $true $null

# 4) Data rows: restyle L (it previously had its own dedicated style) to match
#    the rest of the row, and add the new M (2021) value with the same style.
#    Row 4, 7, 10 are the bold "customFormat" rows; row 4 additionally carries
#    the 0.0 number format.

# Row 4 - bold, 0.0 number format, no border
$ws.Range("L4").Value2 = 10
$ws.Range("M4").Value2 = 10.8
Set-CellLook $ws.Range("L4") $true $false "0.0"
Set-CellLook $ws.Range("M4") $true $false "0.0"

# Row 5 - plain, no border
$ws.Range("L5").Value2 = 6.4
$ws.Range("M5").Value2 = 5.2
Set-CellLook $ws.Range("L5") $false $false $null
Set-CellLook $ws.Range("M5") $false $false $null

# Row 6 - plain, no border
$ws.Range("L6").Value2 = 13.5
$ws.Range("M6").Value2 = 16.2
Set-CellLook $ws.Range("L6") $false $false $null
Set-CellLook $ws.Range("M6") $false $false $null

# Row 7 - bold, no border
$ws.Range("L7").Value2 = 24.3
$ws.Range("M7").Value2 = 24.2
Set-CellLook $ws.Range("L7") $true $false $null
Set-CellLook $ws.Range("M7") $true $false $null

# Row 8 - plain, no border
$ws.Range("L8").Value2 = 27.8
$ws.Range("M8").Value2 = 27.6
Set-CellLook $ws.Range("L8") $false $false $null
Set-CellLook $ws.Range("M8") $false $false $null

# Row 9 - plain, no border
$ws.Range("L9").Value2 = 20.9
$ws.Range("M9").Value2 = 20.9
Set-CellLook $ws.Range("L9") $false $false $null
Set-CellLook $ws.Range("M9") $false $false $null

# Row 10 - bold, no border
$ws.Range("L10").Value2 = 26.7
$ws.Range("M10").Value2 = 28.5
Set-CellLook $ws.Range("L10") $true $false $null
Set-CellLook $ws.Range("M10") $true $false $null

# Row 11 - plain, no border
$ws.Range("L11").Value2 = 28.4
$ws.Range("M11").Value2 = 29.7
Set-CellLook $ws.Range("L11") $false $false $null
Set-CellLook $ws.Range("M11") $false $false $null

# Row 12 - plain, bottom border (last row of the table)
$ws.Range("L12").Value2 = 25
$ws.Range("M12").Value2 = 27.5
Set-CellLook $ws.Range("L12") $false $true $null
Set-CellLook $ws.Range("M12") $false $true $null

# 5) Match the saved selection from the authored workbook.
$ws.Range("O2").Select()
